$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (row 1) renames ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- C95: was string ".." -> numeric 0 ---
$ws.Range("C95").Value = 0

# --- Numeric data updates (GDP column C, and Colony column AL) ---
$ws.Range("C2").Value = 2870.311589353206
$ws.Range("C3").Value = 1873.394108966653
$ws.Range("C4").Value = 1460.056109840828
$ws.Range("C5").Value = 2934.187009790061
$ws.Range("C6").Value = 697.6889104500298
$ws.Range("AL6").Value = 1
$ws.Range("C7").Value = 5191.140356354663
$ws.Range("AL7").Value = 1
$ws.Range("C8").Value = 8947.741473873051
$ws.Range("C9").Value = 4729.735976516416
$ws.Range("C10").Value = 4547.50930098406
$ws.Range("C11").Value = 9502.243585046588
$ws.Range("C12").Value = 5741.405300355145
$ws.Range("C13").Value = 3928.450391496945
$ws.Range("C14").Value = 471.181692645893
$ws.Range("C15").Value = 815.8736791314819
$ws.Range("C16").Value = 2983.242707849043
$ws.Range("C17").Value = 2898.942214704482
$ws.Range("C18").Value = 665.6274194933962
$ws.Range("AL18").Value = 1
$ws.Range("C19").Value = 1904.346464968814
$ws.Range("C20").Value = 9271.398233246389
$ws.Range("C21").Value = 1503.870423231357
$ws.Range("C22").Value = 5555.389721901988
$ws.Range("AL22").Value = 1
$ws.Range("C23").Value = 4633.590358399045
$ws.Range("C24").Value = 4355.934938677345
$ws.Range("C25").Value = 5082.354756663512
$ws.Range("C26").Value = 4094.350334420203
$ws.Range("C27").Value = 492.3430015592067
$ws.Range("C28").Value = 864.5379000312432
$ws.Range("C29").Value = 3083.80337578809
$ws.Range("C30").Value = 2965.153206179127
$ws.Range("C31").Value = 691.8942672110555
$ws.Range("AL31").Value = 1
$ws.Range("C32").Value = 1939.33862702996
$ws.Range("C33").Value = 9477.887185090232
$ws.Range("C34").Value = 5660.517066940175
$ws.Range("AL34").Value = 1
$ws.Range("C35").Value = 4921.848409120176
$ws.Range("C36").Value = 4479.398934239905
$ws.Range("C37").Value = 5360.226632400601
$ws.Range("C38").Value = 14179.19231490798
$ws.Range("C39").Value = 513.7390871590731
$ws.Range("C40").Value = 869.6014949562591
$ws.Range("C41").Value = 3156.723844635973
$ws.Range("C42").Value = 701.4459636783288
$ws.Range("AL42").Value = 1
$ws.Range("C43").Value = 9690.869064532331
$ws.Range("C44").Value = 1982.009737844954
$ws.Range("C45").Value = 5745.422744292303
$ws.Range("AL45").Value = 1
$ws.Range("C46").Value = 5122.180090208862
$ws.Range("C47").Value = 4394.543881413723
$ws.Range("C48").Value = 5642.578115155247
$ws.Range("C49").Value = 11745.7759262897
$ws.Range("C50").Value = 534.5063430177229
$ws.Range("C51").Value = 872.1235974568563
$ws.Range("C52").Value = 3212.740625904757
$ws.Range("C53").Value = 720.7128711178943
$ws.Range("AL53").Value = 1
$ws.Range("C54").Value = 9693.722968944676
$ws.Range("C55").Value = 2000.792448761861
$ws.Range("C56").Value = 5955.175904294275
$ws.Range("AL56").Value = 1
$ws.Range("C57").Value = 5295.682695961288
$ws.Range("C58").Value = 4699.493713911862
$ws.Range("C59").Value = 5919.20956823756
$ws.Range("C60").Value = 11993.48398487312
$ws.Range("C61").Value = 11951.20944634967
$ws.Range("C62").Value = 5412.131646018807
$ws.Range("C63").Value = 3252.634165082374
$ws.Range("C64").Value = 730.3063521039821
$ws.Range("AL64").Value = 1
$ws.Range("C65").Value = 2025.814194788851
$ws.Range("C66").Value = 9839.050190896
$ws.Range("C67").Value = 558.2093442539386
$ws.Range("C68").Value = 4861.287098802361
$ws.Range("C69").Value = 5996.49696468919
$ws.Range("C70").Value = 6301.696269820412
$ws.Range("AL70").Value = 1
$ws.Range("C71").Value = 886.4370030633224
$ws.Range("C72").Value = 11431.15448084494
$ws.Range("C73").Value = 5330.539154475424
$ws.Range("C74").Value = 3314.741082534716
$ws.Range("C75").Value = 729.1196658666737
$ws.Range("AL75").Value = 1
$ws.Range("C76").Value = 2067.29003376698
$ws.Range("C77").Value = 10037.20149040966
$ws.Range("C78").Value = 579.0880693780265
$ws.Range("C79").Value = 4944.191641077407
$ws.Range("C80").Value = 6114.227214287786
$ws.Range("C81").Value = 6661.86504232374
$ws.Range("AL81").Value = 1
$ws.Range("C82").Value = 900.3889853519216
$ws.Range("C83").Value = 1134.924536209078
$ws.Range("C84").Value = 10965.97426143915
$ws.Range("C85").Value = 5176.058803160127
$ws.Range("C86").Value = 3382.563653843273
$ws.Range("C87").Value = 729.8559996981501
$ws.Range("AL87").Value = 1
$ws.Range("C88").Value = 2111.193164269742
$ws.Range("C89").Value = 6411.986543373589
$ws.Range("C90").Value = 10205.79575322194
$ws.Range("C91").Value = 584.2111078769213
$ws.Range("C92").Value = 5089.61202008711
$ws.Range("C93").Value = 6262.368904654469
$ws.Range("C94").Value = 7026.178156858586
$ws.Range("AL94").Value = 1
$ws.Range("C96").Value = 909.5979669529498
